$p = $ppt.ActivePresentation

# EMU -> points helper. A tiny +0.5 EMU nudge compensates for the runtime
# truncating (rather than rounding) the point value back to EMU on save,
# so the persisted EMU lands exactly on the intended target value.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

# Look up a shape in a Shapes collection by its stable cNvPr Id rather than
# a positional index.
function Get-ShapeById($shapes, $id) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $candidate = $shapes.Item($k)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Refresh the "datetimeFigureOut" date footer field everywhere it is
#    defined: once per slide layout, and once on the slide master.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "01/12/2025"
        }
    }
}

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $shp = $master.Shapes.Item($si)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "01/12/2025"
    }
}

# ---------------------------------------------------------------------------
# 2. Update the MicroUI diagram on the single content slide: rewording of
#    the "user / platform" legend boxes (drop the 'platform' terminology)
#    plus the accompanying resize/reposition of their text boxes and the
#    small coloured marker swatch next to the "Embedded" label.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "TextBox 76" (id 77) - "Provided by user" -> "Provided by user(s)"
$box1 = Get-ShapeById $s.Shapes 77
$box1.TextFrame.TextRange.Text = "Provided by user(s)"
$box1.Width = EmuToPt 1422327

# "TextBox 77" (id 78) - "Provided by platform" -> "Provided by VEE Port"
$box2 = Get-ShapeById $s.Shapes 78
$box2.TextFrame.TextRange.Text = "Provided by VEE Port"
$box2.Left = EmuToPt 3958094
$box2.Width = EmuToPt 1564166

# "TextBox 78" (id 79) - "Platform" -> "Embedded"
$box3 = Get-ShapeById $s.Shapes 79
$box3.TextFrame.TextRange.Text = "Embedded"
$box3.Left = EmuToPt 3378217
$box3.Top = EmuToPt 1988920
$box3.Width = EmuToPt 958192

# "Rounded Rectangle 164" (id 87, orange marker) shifts right along with
# the relabelled legend text next to it.
$marker = Get-ShapeById $s.Shapes 87
$marker.Left = EmuToPt 3777072
